# Apply the "MCP6004 -> MCP6074" opamp change to the electronic load BOM workbook.

$wb  = $excel.ActiveWorkbook
$wsBOM = $wb.Worksheets.Item("electronicloadBOM")
$wsRev = $wb.Worksheets.Item("Revision")

# --- 1. Update the opamp row (row 12) on the BOM sheet ---------------------
# Value column (B12): MCP6004-I/SL -> MCP6074-E/SL
$wsBOM.Range("B12").Value = "MCP6074-E/SL"

# Description column (E12): "1 MHz, Low-Power Op Amp" -> "Op Amp"
$wsBOM.Range("E12").Value = "Op Amp"

# MPN column (H12): MCP6004-I/SL -> MCP6074-E/SL, adopt the highlighted
# "Neutral" style that the Value cell (B12) already carries.
$wsBOM.Range("B12").Copy() | Out-Null
$wsBOM.Range("H12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsBOM.Range("H12").Value = "MCP6074-E/SL"

# Ebay column (I12): the old ebay listing link is gone now, clear it but
# keep its (hyperlink-blue) cell style.
$wsBOM.Range("I12").ClearContents() | Out-Null

$excel.CutCopyMode = 0

# --- 2. Append the new revision entry on the Revision sheet -----------------
$wsRev.Range("B7").Value = 2.04
$wsRev.Range("C7").Value = "Changed opamp MCP4006 with MCP6074-E/SL"

# --- 3. Update selections / active sheet to match the saved workbook state --
$wsBOM.Range("B12").Select() | Out-Null
$wsRev.Range("C10").Select() | Out-Null
$wsRev.Activate() | Out-Null
